$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new "Save" header in H1, reusing the same formatting (bold font,
# border, centered alignment) already used by the other header cells,
# by copying the direct formatting from G1 (the "sum" header) onto H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add new data cell H2 with value 0 (default/unformatted like the other
# numeric cells in row 2).
$ws.Range("H2").Value = 0
